{"js": "const sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\nconst sec = sections.items[0];\nconst hdr = sec.getHeader(\"Primary\");\nlet out = [];\nout.push(\"header proto: \" + JSON.stringify(Object.getOwnPropertyNames(Object.getPrototypeOf(hdr))));\nreturn out.join(\"\\n\");\n", "ps1": "$d = $word.ActiveDocument\n$sec = $d.Sections.Item(1)\n$h1 = $sec.Headers.Item(1)\n"}
